$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.985.55'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.17%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.700.74'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.88%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.00'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.02%  '

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.16%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3978'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.55%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.15%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.470'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '53.13'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.49%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.003'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.27%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08808'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.39%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '25.97'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.79%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.463'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.20%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001352'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.70%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.957'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.32%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.713.56'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.91%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.82'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.26%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07193'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.67%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.67'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.44%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.319'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.46%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.36'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.31%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.974.08'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.15%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.380'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.91%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.949'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.27%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.69'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +4.85%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.206'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +16.24%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '162.81'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.69%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '150.84'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +9.17%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.359'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.04%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.637'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +27.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.900.07'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.83%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08552'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.90%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.03151'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.21%  '

$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.038'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.51%  '

$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.156'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.20%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2865'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.89%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09580'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.88'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.63%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8262'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.75%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.02'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.83%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.483'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.95%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.13'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.50%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.688'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.52%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7392'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.09%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.245'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.21%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.392'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.54%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08803'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +8.75%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.18%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '139.25'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.01%  '
